$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Remove the green fill highlight from row 7 (A7:F7), keeping left alignment on A7
$ws.Range("A7").Style = "Normal"
$ws.Range("A7").HorizontalAlignment = -4131  # xlLeft
$ws.Range("A7").NumberFormat = "@"

$ws.Range("B7:F7").Style = "Normal"

# Add new data in row 18: C18=5, D18=5, with same style as B18 (green-fill-less, same as B7 style now)
$ws.Range("C18").Value = 5
$ws.Range("D18").Value = 5
$ws.Range("A18:F18").Style = "Normal"

# Apply style copy: A18 should look like A1-ish (style 6 = green fill, no format/align) and B18:F18 like style1 (green fill no align)
# But per diff, style 6 is (green fill, numFmt49, left align) and style1 is (green fill no align)
# Use Format Painter approach: copy format from A7(old) - but A7 already changed. Instead set directly.

$ws.Range("A18").Interior.ColorIndex = 10  # placeholder, will refine
